# Apply "typy edycji" update to Rodzaje_dokumentow_Eteczka.xlsx
# - column D switches from the placeholder shared text "aa" to a numeric
#   "edit type" (1 or 2) per row
# - a new little lookup/header block is written into G1:P1 (field names)
#   and G2:M2 (sample row) next to the existing table
# - view/pane settings follow the freshly edited area

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row values for column D (the new numeric "edit type")
$dValues = @{
    1=1; 2=1; 3=1; 4=1; 5=1; 6=2; 7=1; 8=1; 9=1; 10=1;
    11=1; 12=1; 13=1; 14=1; 15=1; 16=1; 17=1; 18=1; 19=1; 20=1;
    21=1; 22=1; 23=1; 24=1; 25=2; 26=2; 27=2; 28=1; 29=1; 30=1;
    31=1; 32=1; 33=1; 34=1; 35=1; 36=1; 37=1; 38=1; 39=1; 40=1;
    41=1; 42=1; 43=1; 44=1; 45=1; 46=1; 47=1; 48=1; 49=1; 50=1;
    51=2; 52=1; 53=1; 54=1; 55=1; 56=1; 57=1; 58=1; 59=1; 60=1;
    61=1; 62=1; 63=2; 64=2; 65=1; 66=1; 67=1; 68=1; 69=2; 70=1
}

for ($r = 1; $r -le 70; $r++) {
    $ws.Cells.Item($r, 4).Value = $dValues[$r]
}

# New header row describing the import columns used to reload the table
# into the external system ("slon"). Filled left-to-right starting at the
# "symbol" column, matching how the columns were laid out when editing;
# the new "typ" column (G) and the duplicated sample row (2) are filled
# in afterwards, which is also why they land at the end of the shared
# string table.
$ws.Cells.Item(1, 8).Value  = "symbol"
$ws.Cells.Item(1, 9).Value  = "nazwa"
$ws.Cells.Item(1, 10).Value = "dokwlasny"
$ws.Cells.Item(1, 11).Value = "teczkadzial"
$ws.Cells.Item(1, 12).Value = "firma"
$ws.Cells.Item(1, 13).Value = "pracownik"
$ws.Cells.Item(1, 14).Value = "datadokumentu"
$ws.Cells.Item(1, 15).Value = "datapocz"
$ws.Cells.Item(1, 16).Value = "datakoniec"

# Sample/example row right underneath
$ws.Cells.Item(2, 7).Value  = "typ 1"
$ws.Cells.Item(1, 7).Value  = "typ2"
$ws.Cells.Item(2, 8).Value  = "symbol"
$ws.Cells.Item(2, 9).Value  = "nazwa"
$ws.Cells.Item(2, 10).Value = "dokwlasny"
$ws.Cells.Item(2, 11).Value = "teczkadzial"
$ws.Cells.Item(2, 12).Value = "firma "
$ws.Cells.Item(2, 13).Value = "pracownik"

# Column widths for the newly-used columns
$ws.Range("C:C").ColumnWidth = 4
$ws.Range("D:D").ColumnWidth = 2
$ws.Range("I:I").ColumnWidth = 6.5703125
$ws.Range("J:J").ColumnWidth = 10.42578125
$ws.Range("K:K").ColumnWidth = 10.7109375
$ws.Range("M:M").ColumnWidth = 12.85546875
$ws.Range("N:N").ColumnWidth = 15.28515625

# Refresh the frozen pane / active selection to match the area just edited
$ws.Range("A14").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A15").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("G2:M2").Select()
